# Adds the new "structure" to the training-log sheet: extra header columns
# (N:Q), four data rows (3-6) describing VGG experiment runs, a widened
# merged header (J1:Q1) and a couple of new cell styles (date + right align).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.6
$ws.Columns.Item(8).ColumnWidth = 26.6
$ws.Columns.Item(14).ColumnWidth = 10.42
$ws.Columns.Item(15).ColumnWidth = 13.25
$ws.Columns.Item(16).ColumnWidth = 17.75
$ws.Columns.Item(17).ColumnWidth = 42.42

# ---------------------------------------------------------------------
# Row 1 - extend the merged/highlighted title band to the new columns
# ---------------------------------------------------------------------
$ws.Range("J1").Copy()
$ws.Range("N1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("J1:M1").UnMerge()
$ws.Range("J1:Q1").Merge()

# ---------------------------------------------------------------------
# Row 2 - new header labels
# ---------------------------------------------------------------------
$ws.Range("N2").Value = "Optimizer"
$ws.Range("O2").Value = "Function Loss"
$ws.Range("P2").Value = "Linear layers count"
$ws.Range("Q2").Value = "Augmentations"

$ws.Range("B2").Copy()
$ws.Range("N2:Q2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Data rows 3-6
# ---------------------------------------------------------------------

function Set-Text($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

function Set-RightText($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.HorizontalAlignment = -4152
}

function Set-Date($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.NumberFormat = "d-mmm"
}

# --- Row 3 : VGG v1.1 --------------------------------------------------
Set-Date        "A3" "11/30/2024"
Set-Text        "B3" "VGG v1.1"
Set-Text        "C3" "0.001"
Set-Text        "D3" "ReLU"
Set-Text        "E3" "ReLU"
Set-Text        "F3" "3, 2"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 1
Set-Text        "J3" "48х48"
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 3
Set-Text        "N3" "Adam"
Set-Text        "O3" "CE"
Set-RightText   "P3" "__"
Set-RightText   "Q3" "__"

# --- Row 4 : VGG v1.2 --------------------------------------------------
Set-Date        "A4" "12/1/2024"
Set-Text        "B4" "VGG v1.2"
Set-Text        "C4" "0.001"
Set-Text        "D4" "RELU"
Set-Text        "E4" "ReLU"
Set-Text        "F4" "3, 2"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 1
Set-Text        "J4" "48x48"
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 10
Set-Text        "N4" "Adam"
Set-Text        "O4" "CE"
Set-RightText   "P4" "__"
Set-RightText   "Q4" "__"

# --- Row 5 : VGG v1.3 --------------------------------------------------
Set-Date        "A5" "12/1/2024"
Set-Text        "B5" "VGG v1.3"
Set-Text        "C5" "0.001"
Set-Text        "D5" "ReLU"
Set-Text        "E5" "ReLU"
Set-Text        "F5" "3, 2"
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 1
Set-Text        "J5" "64x64"
$ws.Range("K5").Value = 16
$ws.Range("L5").Value = 10
Set-Text        "N5" "Adam"
Set-Text        "O5" "CE"
Set-RightText   "P5" "__"
Set-RightText   "Q5" "__"

# --- Row 6 : VGG v1.4 (best run) ---------------------------------------
Set-Date        "A6" "12/1/2024"
Set-Text        "B6" "VGG v1.4"
Set-Text        "C6" "0.001"
Set-Text        "D6" "ReLU"
Set-Text        "E6" "ReLU"
Set-Text        "F6" "3, 2"
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 6
$ws.Range("I6").Value = 1
Set-Text        "J6" "64x64"
$ws.Range("K6").Value = 24
Set-RightText   "L6" "35(15-20)"
Set-Text        "N6" "Adam"
Set-Text        "O6" "CE"
$ws.Range("P6").Value = 3
Set-RightText   "Q6" "RHF(0.5), RR(10), RVF(0.5), CJ(0.2, 0.2, 0.2, 0.1) "
Set-Text        "R6" "Лучшая"

# ---------------------------------------------------------------------
# Selection left as it was at the end of the editing session
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2:Q2").Select()
